$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Occupation"
$ws.Range("C1").Value = "Amount"
$ws.Range("D1").Value = "Duration"
$ws.Range("E1").Value = "Rate"
$ws.Range("F1").Value = "Monthly Payment"
$ws.Range("G1").Value = "Total Payment"

$ws.Range("D3").Select()
